$p = $ppt.ActivePresentation

# 1) Delete slide 5 (sldId 262 / rId6 - the "MCC Platform" diagram slide)
$p.Slides.Item(5).Delete()

# 2) Update the datetimeFigureOut field text from 04.08.2021 to 12.08.2021
#    across the slide master and all slide layouts.
$oldDate = "04.08.2021"
$newDate = "12.08.2021"

$master = $p.SlideMaster
for ($i = 1; $i -le $master.Shapes.Count; $i++) {
    $shape = $master.Shapes.Item($i)
    if ($shape.HasTextFrame) {
        $tr = $shape.TextFrame.TextRange
        if ($tr.Text -eq $oldDate) {
            $tr.Text = $newDate
        }
    }
}

for ($li = 1; $li -le $master.CustomLayouts.Count; $li++) {
    $layout = $master.CustomLayouts.Item($li)
    for ($i = 1; $i -le $layout.Shapes.Count; $i++) {
        $shape = $layout.Shapes.Item($i)
        if ($shape.HasTextFrame) {
            $tr = $shape.TextFrame.TextRange
            if ($tr.Text -eq $oldDate) {
                $tr.Text = $newDate
            }
        }
    }
}

# 3) Move the top-level group "Gruppieren 2" on slide 3
$s3 = $p.Slides.Item(3)
for ($i = 1; $i -le $s3.Shapes.Count; $i++) {
    $shape = $s3.Shapes.Item($i)
    if ($shape.Name -eq "Gruppieren 2") {
        $shape.Left = 1498147 / 12700
        $shape.Top = 742950 / 12700
    }
}
